$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 3552.3845
$ws.Range("I9").Value = 257
$ws.Range("K9").Value = 257
$ws.Range("M9").Value = -88
$ws.Range("H33").Value = 14188207
$ws.Range("I33").Value = 5785165.5
$ws.Range("K33").Value = 5785165.5
$ws.Range("M33").Value = -5784936.5
$ws.Range("H69").Value = 4955.5557
$ws.Range("I69").Value = 1625
$ws.Range("K69").Value = 4875
$ws.Range("M69").Value = -4001
$ws.Range("H70").Value = 1580
$ws.Range("J70").Value = 1816.6666
$ws.Range("L70").Value = 5449.9998
$ws.Range("N70").Value = -5989.9998
$ws.Range("H72").Value = 4955.5557
$ws.Range("I72").Value = 1625
$ws.Range("K72").Value = 14625
$ws.Range("M72").Value = -10257
$ws.Range("H73").Value = 1580
$ws.Range("J73").Value = 1816.6666
$ws.Range("L73").Value = 5449.9998
$ws.Range("N73").Value = -7321.9998
$ws.Range("H74").Value = 4205.533
$ws.Range("I74").Value = 4760.381
$ws.Range("K74").Value = 4760.381
$ws.Range("M74").Value = -3824.381
$ws.Range("H77").Value = 4205.533
$ws.Range("I77").Value = 4760.381
$ws.Range("K77").Value = 23801.905
$ws.Range("M77").Value = -19121.905
$ws.Range("H96").Value = 1174.6086
$ws.Range("J96").Value = 1388.3
$ws.Range("L96").Value = 4164.9
$ws.Range("N96").Value = -6910.9
$ws.Range("H97").Value = 4491
$ws.Range("J97").Value = 4491
$ws.Range("L97").Value = 13473
$ws.Range("N97").Value = -14465
$ws.Range("H132").Value = 2520.1667
$ws.Range("I132").Value = 2467.074
$ws.Range("K132").Value = 7401.222
$ws.Range("M132").Value = -4871.222
$ws.Range("H135").Value = 2181.0908
$ws.Range("I135").Value = 1874
$ws.Range("K135").Value = 16866
$ws.Range("M135").Value = -14331

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3952.8
$ws.Range("I2").Value = 4292.357
$ws.Range("K2").Value = 4292.357
$ws.Range("M2").Value = -4179.357
$ws.Range("H45").Value = 62684.65
$ws.Range("I45").Value = 114915
$ws.Range("J45").Value = 3925.5
$ws.Range("K45").Value = 114915
$ws.Range("L45").Value = 3925.5
$ws.Range("M45").Value = -114538
$ws.Range("N45").Value = -4679.5
$ws.Range("H74").Value = 826293.0600000001
$ws.Range("I74").Value = 1481.3125
$ws.Range("K74").Value = 1481.3125
$ws.Range("M74").Value = -607.3125
$ws.Range("H77").Value = 826293.0600000001
$ws.Range("I77").Value = 1481.3125
$ws.Range("K77").Value = 7406.5625
$ws.Range("M77").Value = -3038.5625
$ws.Range("H110").Value = 1674.8572
$ws.Range("I110").Value = 649.8461
$ws.Range("K110").Value = 649.8461
$ws.Range("M110").Value = 1395.1539
$ws.Range("H116").Value = 3952.8
$ws.Range("I116").Value = 4292.357
$ws.Range("K116").Value = 4292.357
$ws.Range("M116").Value = -1998.357

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H2").Value = 23332.666
$ws.Range("J2").Value = 9999
$ws.Range("L2").Value = 9999
$ws.Range("N2").Value = -10225
$ws.Range("H3").Value = 3952.8
$ws.Range("I3").Value = 4292.357
$ws.Range("K3").Value = 4292.357
$ws.Range("M3").Value = -4178.357
$ws.Range("H76").Value = 28648.5
$ws.Range("J76").Value = 28648.5
$ws.Range("L76").Value = 28648.5
$ws.Range("N76").Value = -29278.5
$ws.Range("H79").Value = 28648.5
$ws.Range("J79").Value = 28648.5
$ws.Range("L79").Value = 28648.5
$ws.Range("N79").Value = -30832.5
$ws.Range("I80").Value = 166666830
$ws.Range("J80").Value = 129.625
$ws.Range("K80").Value = 166666830
$ws.Range("L80").Value = 129.625
$ws.Range("M80").Value = -166665832
$ws.Range("N80").Value = -2125.625
$ws.Range("I83").Value = 166666830
$ws.Range("J83").Value = 129.625
$ws.Range("K83").Value = 833334150
$ws.Range("L83").Value = 648.125
$ws.Range("M83").Value = -833329158
$ws.Range("N83").Value = -10632.125
$ws.Range("H99").Value = 5165.9653
$ws.Range("I99").Value = 8356.538
$ws.Range("K99").Value = 8356.538
$ws.Range("M99").Value = -6858.538

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2945.6462
$ws.Range("I31").Value = 2306.625
$ws.Range("J31").Value = 3035.3333
$ws.Range("K31").Value = 2306.625
$ws.Range("L31").Value = 3035.3333
$ws.Range("M31").Value = -2011.625
$ws.Range("N31").Value = -3625.3333
$ws.Range("H34").Value = 2945.6462
$ws.Range("I34").Value = 2306.625
$ws.Range("J34").Value = 3035.3333
$ws.Range("K34").Value = 2306.625
$ws.Range("L34").Value = 3035.3333
$ws.Range("M34").Value = -2104.625
$ws.Range("N34").Value = -3439.3333
$ws.Range("H38").Value = 2632.087
$ws.Range("I38").Value = 1974.3334
$ws.Range("J38").Value = 5000
$ws.Range("K38").Value = 1974.3334
$ws.Range("L38").Value = 5000
$ws.Range("M38").Value = -1597.3334
$ws.Range("N38").Value = -5754
$ws.Range("H46").Value = 2632.087
$ws.Range("I46").Value = 1974.3334
$ws.Range("J46").Value = 5000
$ws.Range("K46").Value = 1974.3334
$ws.Range("L46").Value = 5000
$ws.Range("M46").Value = -1763.3334
$ws.Range("N46").Value = -5422
$ws.Range("H58").Value = 2938.68
$ws.Range("I58").Value = 2535.1538
$ws.Range("K58").Value = 2535.1538
$ws.Range("M58").Value = -2332.1538
$ws.Range("H132").Value = 13336063
$ws.Range("I132").Value = 2190.8823
$ws.Range("J132").Value = 41670544
$ws.Range("K132").Value = 6572.646900000001
$ws.Range("L132").Value = 125011632
$ws.Range("M132").Value = -4042.646900000001
$ws.Range("N132").Value = -125016692
$ws.Range("H136").Value = 2938.68
$ws.Range("I136").Value = 2535.1538
$ws.Range("K136").Value = 7605.4614
$ws.Range("M136").Value = -5055.4614

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 1333
$ws.Range("I9").Value = 0
$ws.Range("K9").Value = 0
$ws.Range("M9").ClearContents()
$ws.Range("H13").Value = 1050
$ws.Range("I13").Value = 100
$ws.Range("J13").Value = 2000
$ws.Range("K13").Value = 300
$ws.Range("L13").Value = 6000
$ws.Range("M13").Value = -132
$ws.Range("N13").Value = -6336
$ws.Range("H98").Value = 452.7143
$ws.Range("I98").Value = 296.4
$ws.Range("K98").Value = 889.1999999999999
$ws.Range("M98").Value = 608.8000000000001
$ws.Range("H113").Value = 1038.0714
$ws.Range("I113").Value = 1159.75
$ws.Range("J113").Value = 989.4
$ws.Range("K113").Value = 3479.25
$ws.Range("L113").Value = 2968.2
$ws.Range("M113").Value = -1309.25
$ws.Range("N113").Value = -7308.2
$ws.Range("H123").Value = 11916.85
$ws.Range("I123").Value = 4481.8887
$ws.Range("J123").Value = 18000
$ws.Range("K123").Value = 13445.6661
$ws.Range("L123").Value = 54000
$ws.Range("M123").Value = -10995.6661
$ws.Range("N123").Value = -58900

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 45456172
$ws.Range("I102").Value = 50001640
$ws.Range("J102").Value = 1500
$ws.Range("K102").Value = 50001640
$ws.Range("L102").Value = 1500
$ws.Range("M102").Value = -50000018
$ws.Range("N102").Value = -4744
$ws.Range("H113").Value = 3192.7778
$ws.Range("I113").Value = 3297.5
$ws.Range("K113").Value = 3297.5
$ws.Range("M113").Value = -1127.5
$ws.Range("H122").Value = 11408.625
$ws.Range("I122").Value = 12499.5
$ws.Range("K122").Value = 37498.5
$ws.Range("M122").Value = -35048.5

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1992.875
$ws.Range("I16").Value = 1993.1428
$ws.Range("K16").Value = 1993.1428
$ws.Range("M16").Value = -1823.1428

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H40").Value = 44999.5
$ws.Range("J40").Value = 69999
$ws.Range("L40").Value = 69999
$ws.Range("N40").Value = -70297
$ws.Range("H46").Value = 125518.29
$ws.Range("J46").Value = 125518.29
$ws.Range("L46").Value = 125518.29
$ws.Range("N46").Value = -125980.29
$ws.Range("H81").Value = 5040714.5
$ws.Range("I81").Value = 7560249
$ws.Range("K81").Value = 15120498
$ws.Range("M81").Value = -15119437
$ws.Range("H84").Value = 5040714.5
$ws.Range("I84").Value = 7560249
$ws.Range("K84").Value = 75602490
$ws.Range("M84").Value = -75597186
$ws.Range("H96").Value = 1959.4615
$ws.Range("J96").Value = 2052.6667
$ws.Range("L96").Value = 2052.6667
$ws.Range("N96").Value = -4798.6667
$ws.Range("H107").Value = 166667310
$ws.Range("I107").Value = 868
$ws.Range("J107").Value = 500000200
$ws.Range("K107").Value = 2604
$ws.Range("L107").Value = 1500000600
$ws.Range("M107").Value = -684
$ws.Range("N107").Value = -1500004440
$ws.Range("H126").Value = 4839.1
$ws.Range("I126").Value = 4785.4287
$ws.Range("J126").Value = 4964.3335
$ws.Range("K126").Value = 14356.2861
$ws.Range("L126").Value = 14893.0005
$ws.Range("M126").Value = -11886.2861
$ws.Range("N126").Value = -19833.0005
$ws.Range("H132").Value = 32926.406
$ws.Range("I132").Value = 44920.914
$ws.Range("K132").Value = 134762.742
$ws.Range("M132").Value = -132232.742
$ws.Range("H134").Value = 125518.29
$ws.Range("J134").Value = 125518.29
$ws.Range("L134").Value = 376554.87
$ws.Range("N134").Value = -381624.87
